$wb = $excel.ActiveWorkbook

$wsHistory = $wb.Worksheets.Item("Change History")
$wsSurvey  = $wb.Worksheets.Item("eCL Survey")

# --- eCL Survey sheet: update the "Follow up" text for questions 2 and 3 ---
$wsSurvey.Range("F3").Value = "If yes, how?  If no, what suggestions or recommendations could have made it more useful for you?"
$wsSurvey.Rows.Item(3).RowHeight = 60

# --- Change History sheet: add a new row documenting this change ---
$wsHistory.Cells.Item(6, 2).Value = 42451
$wsHistory.Cells.Item(6, 3).Value = "TFS2249 - eCL CSR Survey text changes (modified follow up for questions 2 and 3"
$wsHistory.Cells.Item(6, 4).Value = 1.01
$wsHistory.Cells.Item(6, 5).Value = "Doug Stearns"
$wsHistory.Rows.Item(6).RowHeight = 25.5

$wsSurvey.Range("F4").Value = "If yes, what specifically?  If no, what could have made it more effective or valuable?"
$wsSurvey.Rows.Item(4).RowHeight = 45

# --- Selection / active sheet bookkeeping to match the author's final view ---
[void]$wsHistory.Range("C7").Select()
[void]$wsSurvey.Activate()
[void]$wsSurvey.Range("F4").Select()
